# Add a new "ID" column (C) to the patient sheet: header "ID" in C1 and a
# numeric patient id in C2, matching the structural edit described by the
# source diff (new 3rd column, dimension A1:C2, selection moves to C2).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header + value for column C
$ws.Range("C1").Value = "ID"
$ws.Range("C2").Value = 118112564852

# Give the new column an explicit width (close to the author's column width)
$ws.Columns.Item(3).ColumnWidth = 22.1

# Match the post-edit selection (C2 becomes the active cell)
$ws.Range("C2").Select()
